$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2430.75
$ws.Range("I70").Value = 833.8
$ws.Range("J70").Value = 3571.4285
$ws.Range("K70").Value = 2501.4
$ws.Range("L70").Value = 10714.2855
$ws.Range("M70").Value = -2231.4

$ws.Range("H73").Value = 2430.75
$ws.Range("I73").Value = 833.8
$ws.Range("J73").Value = 3571.4285
$ws.Range("K73").Value = 2501.4
$ws.Range("L73").Value = 10714.2855
$ws.Range("M73").Value = -1565.4

$ws.Range("H112").Value = 2230.7144
$ws.Range("I112").Value = 950
$ws.Range("J112").Value = 2444.1667
$ws.Range("K112").Value = 2850
$ws.Range("L112").Value = 7332.500100000001
$ws.Range("M112").Value = -1742
$ws.Range("N112").Value = -9548.500100000001

$ws.Range("H116").Value = 1382.5834
$ws.Range("I116").Value = 1420
$ws.Range("J116").Value = 1355.8572
$ws.Range("K116").Value = 1420
$ws.Range("L116").Value = 1355.8572
$ws.Range("M116").Value = 2022
$ws.Range("N116").Value = -8239.8572

$ws.Range("H136").Value = 60000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 60000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 60000
$ws.Range("N136").Value = -70200

$ws.Range("H137").Value = 24391308
$ws.Range("I137").Value = 1007.8571
$ws.Range("J137").Value = 166668060
$ws.Range("K137").Value = 3023.5713
$ws.Range("L137").Value = 500004180
$ws.Range("M137").Value = -473.5712999999996
$ws.Range("N137").Value = -500009280

$ws.Range("H138").Value = 3095.976
$ws.Range("I138").Value = 2978.682
$ws.Range("J138").Value = 3225
$ws.Range("K138").Value = 8936.045999999998
$ws.Range("L138").Value = 9675
$ws.Range("M138").Value = -3796.045999999998
$ws.Range("N138").Value = -19955

$ws.Range("H139").Value = 40780
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 40780
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 40780
$ws.Range("N139").Value = -51060

$ws.Range("H140").Value = 42399.5
$ws.Range("I140").Value = 35000
$ws.Range("J140").Value = 49799
$ws.Range("K140").Value = 35000
$ws.Range("L140").Value = 49799
$ws.Range("M140").Value = -29820
$ws.Range("N140").Value = -60159

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2991.3333
$ws.Range("I88").Value = 3084.4
$ws.Range("J88").Value = 2875
$ws.Range("K88").Value = 3084.4
$ws.Range("L88").Value = 2875
$ws.Range("M88").Value = -2678.4
$ws.Range("N88").Value = -3687

$ws.Range("H91").Value = 2991.3333
$ws.Range("I91").Value = 3084.4
$ws.Range("J91").Value = 2875
$ws.Range("K91").Value = 3084.4
$ws.Range("L91").Value = 2875
$ws.Range("M91").Value = -1680.4
$ws.Range("N91").Value = -5683

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 10166.667
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 10166.667
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 10166.667
$ws.Range("N49").Value = -10644.667

$ws.Range("H68").Value = 27967
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 27967
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 27967
$ws.Range("N68").Value = -29589

$ws.Range("H71").Value = 27967
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 27967
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 83901
$ws.Range("N71").Value = -92013

$ws.Range("H82").Value = 15531.167
$ws.Range("I82").Value = 8402.333000000001
$ws.Range("J82").Value = 22660
$ws.Range("K82").Value = 8402.333000000001
$ws.Range("L82").Value = 22660
$ws.Range("M82").Value = -8019.333000000001
$ws.Range("N82").Value = -23426

$ws.Range("H85").Value = 15531.167
$ws.Range("I85").Value = 8402.333000000001
$ws.Range("J85").Value = 22660
$ws.Range("K85").Value = 8402.333000000001
$ws.Range("L85").Value = 22660
$ws.Range("M85").Value = -7076.333000000001
$ws.Range("N85").Value = -25312

$ws.Range("H86").Value = 3680.6365
$ws.Range("I86").Value = 3181.6667
$ws.Range("J86").Value = 4279.4
$ws.Range("K86").Value = 3181.6667
$ws.Range("L86").Value = 4279.4
$ws.Range("M86").Value = -2058.6667
$ws.Range("N86").Value = -6525.4

$ws.Range("H89").Value = 3680.6365
$ws.Range("I89").Value = 3181.6667
$ws.Range("J89").Value = 4279.4
$ws.Range("K89").Value = 15908.3335
$ws.Range("L89").Value = 21397
$ws.Range("M89").Value = -10292.3335
$ws.Range("N89").Value = -32629

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2096.5625
$ws.Range("I31").Value = 1521
$ws.Range("J31").Value = 3362.8
$ws.Range("K31").Value = 1521
$ws.Range("L31").Value = 3362.8
$ws.Range("M31").Value = -1226

$ws.Range("H34").Value = 2096.5625
$ws.Range("I34").Value = 1521
$ws.Range("J34").Value = 3362.8
$ws.Range("K34").Value = 1521
$ws.Range("L34").Value = 3362.8
$ws.Range("M34").Value = -1319

$ws.Range("H50").Value = 10494.167
$ws.Range("I50").Value = 8721
$ws.Range("J50").Value = 12267.333
$ws.Range("K50").Value = 8721
$ws.Range("L50").Value = 12267.333
$ws.Range("M50").Value = -8096
$ws.Range("N50").Value = -13517.333

$ws.Range("H51").Value = 12000.3
$ws.Range("I51").Value = 8950
$ws.Range("J51").Value = 12762.875
$ws.Range("K51").Value = 8950
$ws.Range("L51").Value = 12762.875
$ws.Range("M51").Value = -8214
$ws.Range("N51").Value = -14234.875

$ws.Range("H61").Value = 12000.3
$ws.Range("I61").Value = 8950
$ws.Range("J61").Value = 12762.875
$ws.Range("K61").Value = 8950
$ws.Range("L61").Value = 12762.875
$ws.Range("M61").Value = -8602
$ws.Range("N61").Value = -13458.875

$ws.Range("H105").Value = 1060.7142
$ws.Range("I105").Value = 675
$ws.Range("J105").Value = 1575
$ws.Range("K105").Value = 675
$ws.Range("L105").Value = 1575
$ws.Range("M105").Value = 1072
$ws.Range("N105").Value = -5069

$ws.Range("H109").Value = 44500
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 44500
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 44500
$ws.Range("N109").Value = -46580

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 523.53125
$ws.Range("I113").Value = 498.46667
$ws.Range("J113").Value = 545.64703
$ws.Range("K113").Value = 1495.40001
$ws.Range("L113").Value = 1636.94109
$ws.Range("M113").Value = 674.5999899999999
$ws.Range("N113").Value = -5976.94109

$ws.Range("H123").Value = 2437.5
$ws.Range("I123").Value = 1750
$ws.Range("J123").Value = 3125
$ws.Range("K123").Value = 5250
$ws.Range("L123").Value = 9375
$ws.Range("M123").Value = -2800

$ws.Range("H129").Value = 544.44446
$ws.Range("I129").Value = 487.5
$ws.Range("J129").Value = 1000
$ws.Range("K129").Value = 1462.5
$ws.Range("L129").Value = 3000
$ws.Range("M129").Value = 3537.5
$ws.Range("N129").Value = -13000

$ws.Range("H130").Value = 2023.3334
$ws.Range("I130").Value = 1801.6666
$ws.Range("J130").Value = 2466.6667
$ws.Range("K130").Value = 5404.9998
$ws.Range("L130").Value = 7400.000100000001
$ws.Range("M130").Value = -384.9997999999996
$ws.Range("N130").Value = -17440.0001

$ws.Range("H131").Value = 1866.2748
$ws.Range("I131").Value = 5061.8184
$ws.Range("J131").Value = 1426.8875
$ws.Range("K131").Value = 15185.4552
$ws.Range("L131").Value = 4280.6625
$ws.Range("M131").Value = -10145.4552
$ws.Range("N131").Value = -14360.6625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 4252.0835
$ws.Range("I36").Value = 1508.3334
$ws.Range("J36").Value = 5166.6665
$ws.Range("K36").Value = 1508.3334
$ws.Range("L36").Value = 5166.6665
$ws.Range("M36").Value = -1023.3334
$ws.Range("N36").Value = -6136.6665

$ws.Range("H70").Value = 3777.3845
$ws.Range("I70").Value = 3400.6667
$ws.Range("J70").Value = 4625
$ws.Range("K70").Value = 3400.6667
$ws.Range("L70").Value = 4625
$ws.Range("M70").Value = -3130.6667
$ws.Range("N70").Value = -5165

$ws.Range("H73").Value = 3777.3845
$ws.Range("I73").Value = 3400.6667
$ws.Range("J73").Value = 4625
$ws.Range("K73").Value = 3400.6667
$ws.Range("L73").Value = 4625
$ws.Range("M73").Value = -2464.6667
$ws.Range("N73").Value = -6497

$ws.Range("H99").Value = 18131.143
$ws.Range("I99").Value = 13599.8
$ws.Range("J99").Value = 29459.5
$ws.Range("K99").Value = 13599.8
$ws.Range("L99").Value = 29459.5
$ws.Range("M99").Value = -11353.8
$ws.Range("N99").Value = -33951.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 895.61536
$ws.Range("I136").Value = 466.4762
$ws.Range("J136").Value = 2698
$ws.Range("K136").Value = 1399.4286
$ws.Range("L136").Value = 8094
$ws.Range("M136").Value = 1150.5714

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
